# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot run).
#
# The sheet stores Price (col D) and Volume(1h) (col E) as plain TEXT,
# even when a price looks like a decimal number (e.g. "312.24") -- that is
# how the upstream scraper (openpyxl, t="inlineStr") wrote the workbook.
# If we assign such a numeric-looking literal straight to .Value, Excel
# auto-converts it to a real number and silently changes the cell's type.
# To avoid that, numeric-looking Price values are entered with a leading
# apostrophe (Excel's standard "force text" input marker -- the apostrophe
# itself is never stored), and then the cell Style is reset to "Normal" so
# we don't leave a stray quote-prefix format behind. Percent strings in
# column E (e.g. "  -0.33%  ") are padded with spaces and already can't be
# parsed as plain numbers, so they can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.364.63'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '2.367.99'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'" + '312.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("D6").Value = "'" + '107.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.38%  '

$ws.Range("E7").Value = '  -1.25%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -3.42%  '

$ws.Range("D10").Value = "'" + '40.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.92%  '

$ws.Range("D11").Value = "'" + '0.0916'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("E12").Value = '  -2.40%  '

$ws.Range("E13").Value = '  +1.14%  '

$ws.Range("E14").Value = '  -3.90%  '

$ws.Range("D15").Value = '2.728.75'
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("E16").Value = '  -2.79%  '

$ws.Range("D17").Value = '2.365.20'
$ws.Range("E17").Value = '  -1.20%  '

$ws.Range("D18").Value = '45.387.10'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = "'" + '14.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.58%  '

$ws.Range("E20").Value = '  -1.66%  '

$ws.Range("D21").Value = "'" + '7.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.63%  '

$ws.Range("D22").Value = "'" + '73.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.49%  '

$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("D24").Value = "'" + '259.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.23%  '

$ws.Range("D25").Value = "'" + '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = "'" + '11.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.36%  '

$ws.Range("D28").Value = "'" + '7.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.33%  '

$ws.Range("E29").Value = '  -1.74%  '

$ws.Range("D30").Value = "'" + '0.0970'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.14%  '

$ws.Range("D31").Value = "'" + '22.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.68%  '

$ws.Range("D32").Value = "'" + '36.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.01%  '

$ws.Range("D33").Value = "'" + '166.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.03%  '

$ws.Range("D34").Value = "'" + '2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.20%  '

$ws.Range("E35").Value = '  -2.06%  '

$ws.Range("E36").Value = '  +0.57%  '

$ws.Range("D37").Value = "'" + '4.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.22%  '

$ws.Range("D38").Value = "'" + '1.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.86%  '

$ws.Range("D39").Value = "'" + '3.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").Value = "'" + '2.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.49%  '

$ws.Range("D41").Value = "'" + '0.0353'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.44%  '

$ws.Range("D42").Value = "'" + '98.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.64%  '

$ws.Range("D43").Value = "'" + '69.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.46%  '

# Row 44/45: coin ranking reshuffled -- FirstDigitalUSD now ranks above Algorand.
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = "'" + '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'" + '0.225'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.16%  '

$ws.Range("D46").Value = "'" + '12.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.89%  '

$ws.Range("D47").Value = '1.811.17'
$ws.Range("E47").Value = '  +9.33%  '

$ws.Range("D48").Value = "'" + '83.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.35%  '

$ws.Range("D49").Value = "'" + '5.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.72%  '

$ws.Range("D50").Value = "'" + '111.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.99%  '

$ws.Range("D51").Value = "'" + '9.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.18%  '
